$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column K: one more year of data (2021), mirroring column J ---
$ws.Range("K4").Value = 2021
$ws.Range("K5").Value = 7999.5
$ws.Range("K6").Value = $null
$ws.Range("K7").Formula = "=K5-K8"
$ws.Range("K8").Value = 252.9
$ws.Range("K9").Value = $null
$ws.Range("K10").Value = 690.4
$ws.Range("K11").Value = 968.2
$ws.Range("K12").Value = 655
$ws.Range("K13").Value = 691.2
$ws.Range("K14").Value = 1248.5
$ws.Range("K15").Value = 959.1
$ws.Range("K16").Value = 2596.6
$ws.Range("K17").Value = 133.6
$ws.Range("K18").Value = 57

# --- Match column K's cell formatting to column J's per-row styles ---
$ws.Range("K4").Font.Name = "Times New Roman"
$ws.Range("K4").Font.Size = 9
$ws.Range("K4").Font.Bold = $true
$ws.Range("K4").Font.ThemeColor = 1
$ws.Range("K4").HorizontalAlignment = -4152
$ws.Range("K4").VerticalAlignment = -4108
$ws.Range("K4").WrapText = $true
$ws.Range("K4").Borders.Item(8).LineStyle = -4138
$ws.Range("K4").Borders.Item(8).Weight = -4138
$ws.Range("K4").Borders.Item(9).LineStyle = -4138
$ws.Range("K4").Borders.Item(9).Weight = -4138

$ws.Range("K5").Font.Name = "Times New Roman"
$ws.Range("K5").Font.Size = 9
$ws.Range("K5").Font.Bold = $true
$ws.Range("K5").Font.Color = 0
$ws.Range("K5").NumberFormat = "0.0"
$ws.Range("K5").HorizontalAlignment = -4152
$ws.Range("K5").VerticalAlignment = -4108
$ws.Range("K5").WrapText = $true

$ws.Range("K6").Font.Name = "Times New Roman"
$ws.Range("K6").Font.Size = 9
$ws.Range("K6").Font.Bold = $false
$ws.Range("K6").Font.Color = 0
$ws.Range("K6").NumberFormat = "0.0"
$ws.Range("K6").HorizontalAlignment = -4152
$ws.Range("K6").VerticalAlignment = -4108
$ws.Range("K6").WrapText = $true

$ws.Range("K7").Font.Name = "Times New Roman"
$ws.Range("K7").Font.Size = 9
$ws.Range("K7").Font.Bold = $false
$ws.Range("K7").Font.Color = 0
$ws.Range("K7").NumberFormat = "0.0"
$ws.Range("K7").HorizontalAlignment = -4152
$ws.Range("K7").VerticalAlignment = -4108
$ws.Range("K7").WrapText = $true

$ws.Range("K8").Font.Name = "Times New Roman"
$ws.Range("K8").Font.Size = 9
$ws.Range("K8").Font.Bold = $false
$ws.Range("K8").Font.Color = 0
$ws.Range("K8").NumberFormat = "0.0"
$ws.Range("K8").HorizontalAlignment = -4152
$ws.Range("K8").VerticalAlignment = -4108
$ws.Range("K8").WrapText = $true

$ws.Range("K9").Font.Name = "Times New Roman"
$ws.Range("K9").Font.Size = 9
$ws.Range("K9").Font.Bold = $true
$ws.Range("K9").Font.ThemeColor = 1
$ws.Range("K9").NumberFormat = "0.0"
$ws.Range("K9").HorizontalAlignment = -4152
$ws.Range("K9").VerticalAlignment = -4108
$ws.Range("K9").WrapText = $true

$ws.Range("K10").Font.Name = "Times New Roman"
$ws.Range("K10").Font.Size = 9
$ws.Range("K10").Font.Bold = $false
$ws.Range("K10").Font.ThemeColor = 1
$ws.Range("K10").NumberFormat = "0.0"
$ws.Range("K10").HorizontalAlignment = -4152
$ws.Range("K10").VerticalAlignment = -4108
$ws.Range("K10").WrapText = $true

$ws.Range("K11").Font.Name = "Times New Roman"
$ws.Range("K11").Font.Size = 9
$ws.Range("K11").Font.Bold = $false
$ws.Range("K11").Font.ThemeColor = 1
$ws.Range("K11").NumberFormat = "0.0"
$ws.Range("K11").HorizontalAlignment = -4152
$ws.Range("K11").VerticalAlignment = -4108
$ws.Range("K11").WrapText = $true

$ws.Range("K12").Font.Name = "Times New Roman"
$ws.Range("K12").Font.Size = 9
$ws.Range("K12").Font.Bold = $false
$ws.Range("K12").Font.ThemeColor = 1
$ws.Range("K12").NumberFormat = "0.0"
$ws.Range("K12").HorizontalAlignment = -4152
$ws.Range("K12").VerticalAlignment = -4108
$ws.Range("K12").WrapText = $true

$ws.Range("K13").Font.Name = "Times New Roman"
$ws.Range("K13").Font.Size = 9
$ws.Range("K13").Font.Bold = $false
$ws.Range("K13").Font.ThemeColor = 1
$ws.Range("K13").NumberFormat = "0.0"
$ws.Range("K13").HorizontalAlignment = -4152
$ws.Range("K13").VerticalAlignment = -4108
$ws.Range("K13").WrapText = $true

$ws.Range("K14").Font.Name = "Times New Roman"
$ws.Range("K14").Font.Size = 9
$ws.Range("K14").Font.Bold = $false
$ws.Range("K14").Font.ThemeColor = 1
$ws.Range("K14").NumberFormat = "0.0"
$ws.Range("K14").HorizontalAlignment = -4152
$ws.Range("K14").VerticalAlignment = -4108
$ws.Range("K14").WrapText = $true

$ws.Range("K15").Font.Name = "Times New Roman"
$ws.Range("K15").Font.Size = 9
$ws.Range("K15").Font.Bold = $false
$ws.Range("K15").Font.ThemeColor = 1
$ws.Range("K15").NumberFormat = "0.0"
$ws.Range("K15").HorizontalAlignment = -4152
$ws.Range("K15").VerticalAlignment = -4108
$ws.Range("K15").WrapText = $true

$ws.Range("K16").Font.Name = "Times New Roman"
$ws.Range("K16").Font.Size = 9
$ws.Range("K16").Font.Bold = $false
$ws.Range("K16").Font.ThemeColor = 1
$ws.Range("K16").NumberFormat = "0.0"
$ws.Range("K16").HorizontalAlignment = -4152
$ws.Range("K16").VerticalAlignment = -4108
$ws.Range("K16").WrapText = $true

$ws.Range("K17").Font.Name = "Times New Roman"
$ws.Range("K17").Font.Size = 9
$ws.Range("K17").Font.Bold = $false
$ws.Range("K17").Font.ThemeColor = 1
$ws.Range("K17").NumberFormat = "0.0"
$ws.Range("K17").HorizontalAlignment = -4152
$ws.Range("K17").VerticalAlignment = -4108
$ws.Range("K17").WrapText = $true

$ws.Range("K18").Font.Name = "Times New Roman"
$ws.Range("K18").Font.Size = 9
$ws.Range("K18").Font.Bold = $false
$ws.Range("K18").Font.ThemeColor = 1
$ws.Range("K18").NumberFormat = "0.0"
$ws.Range("K18").HorizontalAlignment = -4152
$ws.Range("K18").VerticalAlignment = -4108
$ws.Range("K18").WrapText = $true
$ws.Range("K18").Borders.Item(9).LineStyle = -4138
$ws.Range("K18").Borders.Item(9).Weight = -4138


# Move the active selection the way the author's snapshot has it
$ws.Range("N20").Select()
